# Auto-generated edit script applying value changes from the diff
# (columns D = Price, E = Volume(1h), G = Hora, for rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2" "311.93"
Set-TextValue "E2" "0.72%"
Set-TextValue "G2" "15"
Set-TextValue "D3" "37.90"
Set-TextValue "E3" "-0.38%"
Set-TextValue "G3" "15"
Set-TextValue "E4" "0.98%"
Set-TextValue "G4" "15"
Set-TextValue "D5" "0.07928"
Set-TextValue "E5" "0.44%"
Set-TextValue "G5" "15"
Set-TextValue "D6" "4.408"
Set-TextValue "E6" "1.21%"
Set-TextValue "G6" "15"
Set-TextValue "D7" "1.907"
Set-TextValue "E7" "-3.10%"
Set-TextValue "G7" "15"
Set-TextValue "D8" "8.245"
Set-TextValue "E8" "-0.81%"
Set-TextValue "G8" "15"
Set-TextValue "D9" "2.820"
Set-TextValue "E9" "-7.07%"
Set-TextValue "G9" "15"
Set-TextValue "D10" "0.9272"
Set-TextValue "E10" "-0.39%"
Set-TextValue "G10" "15"
Set-TextValue "D11" "0.1198"
Set-TextValue "E11" "-7.70%"
Set-TextValue "G11" "15"
Set-TextValue "E12" "-1.62%"
Set-TextValue "G12" "15"
Set-TextValue "D13" "0.09218"
Set-TextValue "E13" "4.31%"
Set-TextValue "G13" "15"
Set-TextValue "D14" "0.03367"
Set-TextValue "E14" "-1.39%"
Set-TextValue "G14" "15"
Set-TextValue "D15" "0.09622"
Set-TextValue "E15" "-1.25%"
Set-TextValue "G15" "15"
Set-TextValue "D16" "0.001369"
Set-TextValue "E16" "-1.48%"
Set-TextValue "G16" "15"
Set-TextValue "D17" "0.005924"
Set-TextValue "E17" "-0.31%"
Set-TextValue "G17" "15"
Set-TextValue "D18" "3.531"
Set-TextValue "E18" "-1.69%"
Set-TextValue "G18" "15"
Set-TextValue "D19" "0.3446"
Set-TextValue "E19" "0.30%"
Set-TextValue "G19" "15"
Set-TextValue "D20" "5.283"
Set-TextValue "E20" "5.84%"
Set-TextValue "G20" "15"
Set-TextValue "D21" "0.1284"
Set-TextValue "E21" "-0.84%"
Set-TextValue "G21" "15"
Set-TextValue "D22" "0.2592"
Set-TextValue "E22" "4.34%"
Set-TextValue "G22" "15"
Set-TextValue "D23" "0.02103"
Set-TextValue "E23" "180.15%"
Set-TextValue "G23" "15"
Set-TextValue "E24" "1.36%"
Set-TextValue "G24" "15"
Set-TextValue "D25" "0.001250"
Set-TextValue "E25" "2.79%"
Set-TextValue "G25" "15"
Set-TextValue "D26" "0.004282"
Set-TextValue "E26" "-7.33%"
Set-TextValue "G26" "15"
Set-TextValue "D27" "0.0001301"
Set-TextValue "E27" "-63.78%"
Set-TextValue "G27" "15"
Set-TextValue "G28" "15"
Set-TextValue "G29" "15"
Set-TextValue "G30" "15"
Set-TextValue "G31" "15"
Set-TextValue "G32" "15"
Set-TextValue "G33" "15"
Set-TextValue "G34" "15"
Set-TextValue "G35" "15"
Set-TextValue "G36" "15"
Set-TextValue "G37" "15"
Set-TextValue "G38" "15"
Set-TextValue "D39" "0.02111"
Set-TextValue "E39" "-8.33%"
Set-TextValue "G39" "15"
Set-TextValue "D40" "0.05092"
Set-TextValue "E40" "0.69%"
Set-TextValue "G40" "15"
Set-TextValue "D41" "0.007609"
Set-TextValue "E41" "1.47%"
Set-TextValue "G41" "15"
Set-TextValue "D42" "0.009140"
Set-TextValue "E42" "-7.64%"
Set-TextValue "G42" "15"
Set-TextValue "E43" "-0.42%"
Set-TextValue "G43" "15"
Set-TextValue "D44" "0.002102"
Set-TextValue "E44" "4.23%"
Set-TextValue "G44" "15"
Set-TextValue "D45" "0.008658"
Set-TextValue "E45" "-1.23%"
Set-TextValue "G45" "15"
Set-TextValue "D46" "0.00006686"
Set-TextValue "E46" "2.03%"
Set-TextValue "G46" "15"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "E47" "0.10%"
Set-TextValue "G47" "15"
Set-TextValue "D48" "0.002898"
Set-TextValue "E48" "-3.17%"
Set-TextValue "G48" "15"
Set-TextValue "E49" "-0.08%"
Set-TextValue "G49" "15"
Set-TextValue "D50" "0.00002099"
Set-TextValue "E50" "0.10%"
Set-TextValue "G50" "15"
Set-TextValue "D51" "0.0001999"
Set-TextValue "E51" "0.10%"
Set-TextValue "G51" "15"
